$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows to append, describing additional feed entries scraped by the workflow.
$newRows = @(
    @{
        Link    = "https://www.genomeweb.com/cancer/ogt-assay-nabs-fda-authorization-cdx-syndaxs-revuforj-acute-leukemia"
        Keyword = "CDx"
        Title   = "OGT Assay Nabs FDA Authorization as CDx for Syndax's Revuforj in Acute Leukemia"
    },
    @{
        Link    = "https://www.360dx.com/cancer/agilent-technologies-lunit-partner-cancer-cdx-development"
        Keyword = "CDx, companion diagnostics"
        Title   = "Agilent Technologies, Lunit Partner for Cancer CDx Development"
    },
    @{
        Link    = "https://www.360dx.com/cancer/ogt-assay-nabs-fda-authorization-cdx-syndaxs-revuforj-acute-leukemia"
        Keyword = "CDx"
        Title   = "OGT Assay Nabs FDA Authorization as CDx for Syndax's Revuforj in Acute Leukemia"
    }
)

$startRow = 47
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $row = $startRow + $i
    $entry = $newRows[$i]

    # Copy the formatting (incl. the Hyperlink cell style) from the last existing
    # data row so the new row matches the look of the rest of the table.
    $ws.Range("A46").Copy($ws.Range("A" + $row))

    $ws.Range("A" + $row).Value = $entry.Link
    $ws.Range("B" + $row).Value = $entry.Keyword
    $ws.Range("C" + $row).Value = $entry.Title

    $ws.Hyperlinks.Add($ws.Range("A" + $row), $entry.Link)
    $ws.Range("A" + $row).Style = "Hyperlink"
}
